# Rename the worksheet "Sheet1" to "Balaji" (Add files via upload).
$wb = $excel.ActiveWorkbook

$target = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Sheet1") {
        $target = $sheet
        break
    }
}

if ($target -ne $null) {
    $target.Name = "Balaji"
}
